$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append new key/value rows for archaea motility flagella data
$ws.Range("A18").Value = "motilityFlagellaMonotrichous"
$ws.Range("B18").Value = "Monotrichous Flagella"

$ws.Range("A19").Value = "motilityFlagellaLophotrichous"
$ws.Range("B19").Value = "Lophotrichous Flagella"

$ws.Range("A20").Value = "motilityFlagellaPeritrichous"
$ws.Range("B20").Value = "Peritrichous Flagella"

$ws.Range("A21").Value = "motilityFlagellaAmphitrichous"
$ws.Range("B21").Value = "Amphitrichous Flagella"

# Match the selection state recorded in the workbook after the edit
$ws.Range("B18").Select()
